# Generate Report for Handoff
#
# Updates the localization-status report:
#   - On the "zh-cn" sheet, rows 7,8,9,10,11,13 now have a Priority of "ht"
#     and a refreshed "Latest Handoff Datetime".
#   - On the "de-de" sheet, the same rows get the same Priority update and
#     their own refreshed "Latest Handoff Datetime".
#   - The "Overview" sheet's "Latest HO Xliff Generate Date" for those same
#     rows is refreshed to match the new de-de handoff datetime (it shared
#     the same value as the de-de sheet before this edit).

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 10, 11, 13)

# zh-cn sheet
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-12 08:24:44"
}

# de-de sheet
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-12 08:24:51"
}

# Overview sheet
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-12 08:24:51"
}
